# EVR_YR_FIN.xlsx update - "Doing Updates for Financials"
#
# A new reporting period is inserted as column D on the sole worksheet,
# pushing all previously existing periods one column to the right
# (old D -> E, old E -> F, ... old K -> L). The new column D is then
# populated with the newest period's figures. A handful of
# derived/summary rows also get their first couple of shifted-in cells
# recomputed, because the new period changes the underlying totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; this shifts all existing data
# (columns D:K) one column to the right (to E:L), matching the
# worksheet's new A5:L102 dimension.
$ws.Columns("D").Insert()

# ---- Carry over formatting for the new column from its neighbour ----
# (column E, which now holds what used to be column D) so the new
# period's cells keep the same number formats / fonts as the rest of
# the table, instead of Excel's generic default style.
$ws.Range("E8:E35").Copy()
$ws.Range("D8:D35").PasteSpecial(-4122)

$ws.Range("E39:E77").Copy()
$ws.Range("D39:D77").PasteSpecial(-4122)

$ws.Range("E81:E102").Copy()
$ws.Range("D81:D102").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E80").Copy()
$ws.Range("D80").PasteSpecial(-4122)

# ---- Populate the new column D with the newest reporting period ----

# Period Ending (date) headers
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# Income Statement section
$ws.Range("D8").Value = 2082500
$ws.Range("D9").Value = 17800
$ws.Range("D10").Value = 2064700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 5000
$ws.Range("D15").Value = 27100
$ws.Range("D17").Value = 1540400
$ws.Range("D18").Value = 542100
$ws.Range("D20").Value = 9300
$ws.Range("D21").Value = 580700
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 551400
$ws.Range("D24").Value = 108100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 443300
$ws.Range("D27").Value = 377600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -400
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -9300
$ws.Range("D33").Value = 377200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 377200

# Balance Sheet section
$ws.Range("D41").Value = 790600
$ws.Range("D42").Value = 307300
$ws.Range("D43").Value = 334800
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 48900
$ws.Range("D46").Value = 1481600
$ws.Range("D47").Value = 91200
$ws.Range("D48").Value = 81100
$ws.Range("D49").Value = 141800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 330000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2125700
$ws.Range("D57").Value = 37900
$ws.Range("D58").Value = 25100
$ws.Range("D59").Value = 686700
$ws.Range("D60").Value = 749700
$ws.Range("D61").Value = 168600
$ws.Range("D62").Value = 199400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1367500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 364900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 758100
$ws.Range("D77").Value = 0

# Cash Flow section
$ws.Range("D81").Value = 377200
$ws.Range("D83").Value = 29400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 849600
$ws.Range("D91").Value = -33300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -212600
$ws.Range("D96").Value = -77300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -452900
$ws.Range("D101").Value = -1400
$ws.Range("D102").Value = 182700

# A few derived rows also have their next one/two shifted-in cells
# recomputed (not simple carry-overs), because totals changed once the
# new period's figures were added.
$ws.Range("E89").Value = 507200
$ws.Range("F89").Value = 421900
$ws.Range("E94").Value = -54600
$ws.Range("F94").Value = -46200
$ws.Range("E102").Value = 41700
$ws.Range("F102").Value = 112400
